$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mst1"
$ws.Range("C2").Value = "Mst1r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3237273333333333
$ws.Range("H2").Value = 0.9711819999999999
$ws.Range("I2").Value = 0.1800074584373758
$ws.Range("J2").Value = 0.1800074584373758
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2044843333333334
$ws.Range("N2").Value = 0.613453
$ws.Range("O2").Value = 0.03437184025187553
$ws.Range("P2").Value = 0.03437184025187552
$ws.Range("Q2").Value = 0.06619716793844445
$ws.Range("R2").Value = 0.5957745114459999
$ws.Range("S2").Value = 0.006187187605555604
$ws.Range("T2").Value = 0.006187187605555602

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mst1"
$ws.Range("C3").Value = "Mst1r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3237273333333333
$ws.Range("H3").Value = 0.9711819999999999
$ws.Range("I3").Value = 0.1800074584373758
$ws.Range("J3").Value = 0.1800074584373758
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.171219666666667
$ws.Range("N3").Value = 15.513659
$ws.Range("O3").Value = 0.8692320501653281
$ws.Range("P3").Value = 0.869232050165328
$ws.Range("Q3").Value = 1.674065152770889
$ws.Range("R3").Value = 15.066586374938
$ws.Range("S3").Value = 0.1564682521425702
$ws.Range("T3").Value = 0.1564682521425702

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mst1"
$ws.Range("C4").Value = "Mst1r"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3237273333333333
$ws.Range("H4").Value = 0.9711819999999999
$ws.Range("I4").Value = 0.1800074584373758
$ws.Range("J4").Value = 0.1800074584373758
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.573478
$ws.Range("N4").Value = 1.720434
$ws.Range("O4").Value = 0.09639610958279644
$ws.Range("P4").Value = 0.09639610958279642
$ws.Range("Q4").Value = 0.1856505036653333
$ws.Range("R4").Value = 1.670854532988
$ws.Range("S4").Value = 0.01735201868924995
$ws.Range("T4").Value = 0.01735201868924995

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Mst1"
$ws.Range("C5").Value = "Mst1r"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9248883333333334
$ws.Range("H5").Value = 2.774665
$ws.Range("I5").Value = 0.5142809428769699
$ws.Range("J5").Value = 0.5142809428769699
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2044843333333334
$ws.Range("N5").Value = 0.613453
$ws.Range("O5").Value = 0.03437184025187553
$ws.Range("P5").Value = 0.03437184025187552
$ws.Range("Q5").Value = 0.1891251742494445
$ws.Range("R5").Value = 1.702126568245
$ws.Range("S5").Value = 0.01767678241315113
$ws.Range("T5").Value = 0.01767678241315113

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mst1"
$ws.Range("C6").Value = "Mst1r"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9248883333333334
$ws.Range("H6").Value = 2.774665
$ws.Range("I6").Value = 0.5142809428769699
$ws.Range("J6").Value = 0.5142809428769699
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.171219666666667
$ws.Range("N6").Value = 15.513659
$ws.Range("O6").Value = 0.8692320501653281
$ws.Range("P6").Value = 0.869232050165328
$ws.Range("Q6").Value = 4.782800738803889
$ws.Range("R6").Value = 43.04520664923501
$ws.Range("S6").Value = 0.4470294783379066
$ws.Range("T6").Value = 0.4470294783379065

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mst1"
$ws.Range("C7").Value = "Mst1r"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9248883333333334
$ws.Range("H7").Value = 2.774665
$ws.Range("I7").Value = 0.5142809428769699
$ws.Range("J7").Value = 0.5142809428769699
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.573478
$ws.Range("N7").Value = 1.720434
$ws.Range("O7").Value = 0.09639610958279644
$ws.Range("P7").Value = 0.09639610958279642
$ws.Range("Q7").Value = 0.5304031116233334
$ws.Range("R7").Value = 4.77362800461
$ws.Range("S7").Value = 0.04957468212591227
$ws.Range("T7").Value = 0.04957468212591225

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Mst1"
$ws.Range("C8").Value = "Mst1r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.549795
$ws.Range("H8").Value = 1.649385
$ws.Range("I8").Value = 0.3057115986856543
$ws.Range("J8").Value = 0.3057115986856543
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2044843333333334
$ws.Range("N8").Value = 0.613453
$ws.Range("O8").Value = 0.03437184025187553
$ws.Range("P8").Value = 0.03437184025187552
$ws.Range("Q8").Value = 0.112424464045
$ws.Range("R8").Value = 1.011820176405
$ws.Range("S8").Value = 0.01050787023316879
$ws.Range("T8").Value = 0.01050787023316879

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Mst1"
$ws.Range("C9").Value = "Mst1r"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.549795
$ws.Range("H9").Value = 1.649385
$ws.Range("I9").Value = 0.3057115986856543
$ws.Range("J9").Value = 0.3057115986856543
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.171219666666667
$ws.Range("N9").Value = 15.513659
$ws.Range("O9").Value = 0.8692320501653281
$ws.Range("P9").Value = 0.869232050165328
$ws.Range("Q9").Value = 2.843110716635
$ws.Range("R9").Value = 25.587996449715
$ws.Range("S9").Value = 0.2657343196848513
$ws.Range("T9").Value = 0.2657343196848513

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Mst1"
$ws.Range("C10").Value = "Mst1r"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.549795
$ws.Range("H10").Value = 1.649385
$ws.Range("I10").Value = 0.3057115986856543
$ws.Range("J10").Value = 0.3057115986856543
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.573478
$ws.Range("N10").Value = 1.720434
$ws.Range("O10").Value = 0.09639610958279644
$ws.Range("P10").Value = 0.09639610958279642
$ws.Range("Q10").Value = 0.3152953370100001
$ws.Range("R10").Value = 2.83765803309
$ws.Range("S10").Value = 0.02946940876763422
$ws.Range("T10").Value = 0.02946940876763422
